$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update existing parameter values on row 2 ---
$ws.Range("C2").Value = 3
$ws.Range("F2").Value = 20
$ws.Range("G2").Value = 0.5

# --- Add the three new header labels to the shared-string pool ---
$ws.Range("AB1").Value = "min_n"
$ws.Range("AC1").Value = "max_n"
$ws.Range("AD1").Value = "bucket"
$ws.Range("AB1:AD1").Interior.Color = $ws.Range("A1").Interior.Color

# --- Re-create the header row (A1:M1) starting at column O for the new plot/table ---
$headers = @("vector_size","window_size","min_count","sg","hs","negative","ns_exponent","seed","batchs_words","shrink_windows","alpha","workers","epochs")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $col = 15 + $i  # column O = 15
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
    $cell.Interior.Color = $ws.Range("A1").Interior.Color
}

$ws.Range("K24").Select()
